# Finalized the excel importer
# - Adds a new "Aircraft" parameter block (Aircraft / M cruise / EOC) above
#   the existing "Fuselage" block.
# - Updates tailSlenderness (2 -> 2.5) and tailUpAngle (5 -> 10) values.
# - Keeps the sheet's total row count / dimension unchanged (A1:H35) by
#   trimming an equal number of rows from the trailing blank filler block.
# - Moves the active selection to C17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 fresh rows above the "Fuselage" row (old row 2), shifting the
#    Fuselage/Wing/EOF blocks down by 5 rows.
$ws.Rows("2:6").Insert()

# 2) The lone blank row that used to separate "maTechnology" from
#    "wingPosition" (old row 14) is dropped entirely in the target layout -
#    after the insert above it now sits at row 19.
$ws.Rows("19:19").Delete()

# 3) Remove 4 more rows from the trailing blank filler block so the sheet
#    keeps its original dimensions (A1:H35) - 5 inserted, 1+4 removed.
$ws.Rows("36:39").Delete()

# 4) The Insert() operation copies the style (incl. number format) of the
#    row it pushed down onto the new blank rows 2-6 (col C). Strip that
#    back out wherever the target layout wants a truly empty cell.
$ws.Range("C2").Clear()
$ws.Range("C3").Clear()
$ws.Range("C5").Clear()
$ws.Range("C6").Clear()
$ws.Range("C4").Clear()

# 5) Populate the new "Aircraft" parameter block.
$ws.Range("A3").Value = "Aircraft"
$ws.Range("B4").Value = "M cruise"
$ws.Range("C4").Value = 0.77
$ws.Range("E4").Value = "float"
$ws.Range("B5").Value = "EOC"

# 6) Update existing Fuselage values that changed.
$ws.Range("C11").Value = 2.5
$ws.Range("C12").Value = 10

# 7) Restore the active selection to C17 (matches the post-edit selection).
$ws.Range("C17").Select()
